$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The region labels in column A were shuffled/corrected: some regions had the
# wrong (abbreviated / mis-hyphenated) Italian names, and the order of a few
# blocks of 10 rows (one block per region, 2010-2019) was wrong. This
# reassigns the correct official region name to every 10-row block, in the
# corrected order, and fixes the official long-form names for the
# special-statute regions (which also get a distinct, smaller font).

$ws.Range("A2:A11").Value = "Piemonte"
$ws.Range("A22:A31").Value = "Liguria"
$ws.Range("A32:A41").Value = "Lombardia"
$ws.Range("A52:A61").Value = "Veneto"
$ws.Range("A82:A91").Value = "Toscana"
$ws.Range("A92:A101").Value = "Umbria"
$ws.Range("A102:A111").Value = "Marche"
$ws.Range("A112:A121").Value = "Lazio"
$ws.Range("A132:A141").Value = "Molise"
$ws.Range("A142:A151").Value = "Campania"
$ws.Range("A152:A161").Value = "Puglia"
$ws.Range("A162:A171").Value = "Basilicata"
$ws.Range("A172:A181").Value = "Calabria"
$ws.Range("A182:A191").Value = "Sicilia"
$ws.Range("A192:A201").Value = "Sardegna"
$ws.Range("A42:A51").Value = "Trentino-Alto Adige/Südtirol"
$ws.Range("A62:A71").Value = "Friuli-Venezia Giulia"
$ws.Range("A12:A21").Value = "Valle d'Aosta/Vallée d'Aoste"
$ws.Range("A72:A81").Value = "Emilia-Romagna"
$ws.Range("A122:A131").Value = "Abruzzo"

# The four regions with composite / accented official names get a distinct
# (smaller) font, vertically centered.
$specialRanges = "A12:A21", "A42:A51", "A62:A71", "A72:A81"
foreach ($addr in $specialRanges) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Var(--colab-code-font-family)"
    $rng.Font.Size = 8
    $rng.Font.Color = 2171169
    $rng.VerticalAlignment = -4108
}

# Restore the selection / scroll position left by the editor.
$excel.ActiveWindow.ScrollRow = 117
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("H125").Select()

Write-Host "Edit applied"
